$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47 - this shifts existing rows 47..139 down to 48..140
$ws.Rows.Item(47).Insert()

# Populate the new row 47 with the new weekly price entry
$ws.Cells.Item(47, 1).Value = 7
$ws.Cells.Item(47, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(47, 3).Value = "Ñuble"
$ws.Cells.Item(47, 4).Value = 45082
$ws.Cells.Item(47, 5).Value = 16
$ws.Cells.Item(47, 6).Value = 100112031
$ws.Cells.Item(47, 7).Value = "Poroto verde"
$ws.Cells.Item(47, 8).Value = "Magnum"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 15
$ws.Cells.Item(47, 11).Value = 25000
$ws.Cells.Item(47, 12).Value = 25000
$ws.Cells.Item(47, 13).Value = 25000
$ws.Cells.Item(47, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(47, 15).Value = "Perú"
$ws.Cells.Item(47, 16).Value = 1000
$ws.Cells.Item(47, 17).Value = 25
$ws.Cells.Item(47, 18).Value = "Hortaliza"
